$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before F (old F/G shift to G/H, inheriting their
# previous formatting & values unchanged). The new column F inherits
# column E's per-row formatting, which already matches the desired
# result for most rows.
$ws.Columns("F").Insert()

# Header for the new column.
$ws.Range("F1").Value = "REINFORCEMENT "

# Fill in "KEVLAR" for the rows that use it (inherited formatting from
# column E already matches the target for all of these except row 4).
$kevlarRows = @(2,3,4,7,8,12,18,19,24,26,28,29,30,32,33,35,36,37,38,39,41)
foreach ($r in $kevlarRows) {
    $ws.Cells.Item($r, 6).Value = "KEVLAR"
}

# Fill in "NOMEX 4 PLY " for the rows that use it.
$nomexRows = @(5,6,9,10,11,21,22,23,25,31,34)
foreach ($r in $nomexRows) {
    $ws.Cells.Item($r, 6).Value = "NOMEX 4 PLY "
}

# Fill in "---" for the remaining two rows.
$ws.Cells.Item(16, 6).Value = "---"
$ws.Cells.Item(17, 6).Value = "---"

# Fix up formatting so it matches the final layout exactly:
# - Row 4's KEVLAR cell ends up with the plain (non-wrapping) style.
# - Every NOMEX 4 PLY cell uses the plain (non-wrapping) style.
# - The "---" cells (and row 12, which had a quote-prefixed style to begin
#   with) use the quote-prefixed style (like the existing "--" cells).
$ws.Range("A2").Copy() | Out-Null
$plainStyleRows = @(4,5,6,9,10,11,21,22,23,25,31,34)
foreach ($r in $plainStyleRows) {
    $ws.Cells.Item($r, 6).PasteSpecial(-4122)
}

$ws.Range("D7").Copy() | Out-Null
$quotePrefixStyleRows = @(12,16,17)
foreach ($r in $quotePrefixStyleRows) {
    $ws.Cells.Item($r, 6).PasteSpecial(-4122)
}

$excel.CutCopyMode = 0
